$d = $word.ActiveDocument

# Locate the paragraph that begins with "3-" (the SourceTree/clone-url
# instructions) by scanning the paragraph collection for its text, rather
# than a hard-coded index.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("3-")) {
        $targetIndex = $i
        break
    }
}

# Remove that whole paragraph, including its paragraph mark, so its
# content merges into the following paragraph (the one holding the
# "_GoBack" bookmark).
$oldPara = $d.Paragraphs($targetIndex)
$oldPara.Range.Delete()

# The bookmark paragraph now sits at $targetIndex; drop the replacement
# text in front of the bookmark that now starts that paragraph.
$bookmarkPara = $d.Paragraphs($targetIndex)
$bookmarkPara.Range.InsertBefore("It is not yet completed")

# Collapse the two trailing empty paragraphs down to one by removing the
# first of the pair (the true last paragraph mark of the body can't be
# removed, so target the one right after the bookmark paragraph).
$trailingEmpty = $d.Paragraphs($targetIndex + 1)
$trailingEmpty.Range.Delete()
